$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the header labels in B1 and C1 (X_UTM / Y_UTM)
$ws.Range("B1").Value = "Y_UTM"
$ws.Range("C1").Value = "X_UTM"

# Add new header "area" in AF1, matching the header formatting used by
# the rest of row 1 (bold, thin border, centered).
$ws.Range("AF1").Value = "area"
$ws.Range("AE1").Copy()
$ws.Range("AF1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Fill AF2:AF180 with the value 5 for every data row
$ws.Range("AF2:AF180").Value = 5

$wb.Save()
